# Remove the "Level 3", "Level 4", "Level 5", "Level 6 - 7" rubric slides
# (originally slides 2-5) from the Kodu lesson deck. All other slides keep
# their existing order and content.
$p = $ppt.ActivePresentation

# Delete from the highest index down to avoid re-indexing issues.
$p.Slides.Item(5).Delete()   # "Level 6 - 7"
$p.Slides.Item(4).Delete()   # "Level 5"
$p.Slides.Item(3).Delete()   # "Level 4"
$p.Slides.Item(2).Delete()   # "Level 3"
